# Applies the "Se agrego las caracteristicas para obtener y crear nuevos
# procesos de compra" edit: extends the Ubicaciones/Procesos lookup tables
# on the "Datos" sheet with new rows, repoints the workbook's defined
# names at the Excel Tables instead of fixed ranges, and updates the
# "Formato" sheet's data validation / column width / selection to match.

$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("Datos")
$wsFormato = $wb.Worksheets.Item("Formato")

# --- 1. Fill in the new "Datos" sheet values -------------------------------
# Order matters: it reproduces the exact shared-string table layout that
# Excel produced when the author typed these values in (J4 before J3 frees
# up the old slot 31 for the "Redes" text, then J3 gets a fresh slot for the
# corrected "CTT" text, etc).
$wsDatos.Range("J4").Value = "5 - Fisei - Piso1 - Laboratorio Redes"
$wsDatos.Range("J3").Value = "3 - Fisei - Piso1 - Laboratorio CTT"
$wsDatos.Range("J5").Value = "6 - Fisei - Piso2 - Laboratorio Redes 02"
$wsDatos.Range("J6").Value = "7 - Fisei - Piso3 - Administración"
$wsDatos.Range("J7").Value = "8 - Fisei - Piso1 - Coordinación"
$wsDatos.Range("J8").Value = "9 - Fche - Piso1 - Coordinación"
$wsDatos.Range("J9").Value = "10 - Fche - Piso1 - Laboratorio 1"
$wsDatos.Range("J10").Value = "11 - Fcial - Piso1 - Laboratorio Central"
$wsDatos.Range("J11").Value = "12 - Fcial - Piso1 - Laboratorio Bacteriologo"
$wsDatos.Range("J12").Value = "13 - Fcial - Piso2 - Laboratorio 05"

$wsDatos.Range("L4").Value = "2 - PR0002"
$wsDatos.Range("L5").Value = "3 - PR0003"
$wsDatos.Range("L6").Value = "4 - PR0004"
$wsDatos.Range("L7").Value = "5 - PR0005"

# --- 2. Resize the Tables (ListObjects) that back the validation lists -----
$tablaUbicaciones = $wsDatos.ListObjects.Item("Tabla4")
$tablaUbicaciones.Resize($wsDatos.Range("J2:J12"))

$tablaProcesos = $wsDatos.ListObjects.Item("Tabla8")
$tablaProcesos.Resize($wsDatos.Range("L2:L7"))

# --- 3. Point the workbook-level defined names at the Tables ---------------
$wb.Names.Item("Categorias").RefersTo = "=Tabla2[Categorias]"
$wb.Names.Item("MARCAS_INF").RefersTo = "=Tabla6[Marcas_Inf]"
$wb.Names.Item("MARCAS_OFI").RefersTo = "=Tabla5[Marcas_Ofi]"
$wb.Names.Item("Procesos").RefersTo = "=Tabla8[Procesos]"
$wb.Names.Item("TIPOS_INF").RefersTo = "=Tabla3[Tipos_Inf]"
$wb.Names.Item("TIPOS_OFI").RefersTo = "=Tabla7[Tipos_Ofi]"
$wb.Names.Item("UBICACIONES").RefersTo = "=Tabla4[Ubicaciones]"

# --- 4. Widen column J on "Datos" and column D on "Formato" ----------------
$wsDatos.Columns.Item(10).ColumnWidth = 33.67
$wsFormato.Columns.Item(4).ColumnWidth = 33.67

# --- 5. Reorder / retarget the "Formato" sheet data validations ------------
# Move the Ubicaciones validation after the Procesos one, and make it skip
# D4 (that row now uses a value that isn't a plain Ubicaciones pick).
$wsFormato.Range("D3:D20").Validation.Delete()
$wsFormato.Range("D3:D20").Validation.Add(3, 1, 1, "UBICACIONES")
$wsFormato.Range("D3:D20").Validation.ErrorTitle = "Datos erroneos"
$wsFormato.Range("D3:D20").Validation.ErrorMessage = "Elige un dato que este disponible en la lista"
$wsFormato.Range("D4").Validation.Delete()

# --- 6. Restore the selections that were recorded in the saved workbook ----
$wsDatos.Activate()
$wsDatos.Range("L9").Select()
$wsFormato.Activate()
$wsFormato.Range("H12").Select()

Write-Output "Edit applied"
